$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 34704.234
$ws.Range("I51").Value = 6985
$ws.Range("K51").Value = 6985
$ws.Range("M51").Value = -6501
$ws.Range("H53").Value = 400.0909
$ws.Range("I53").Value = 474.42856
$ws.Range("J53").Value = 270
$ws.Range("K53").Value = 474.42856
$ws.Range("L53").Value = 270
$ws.Range("M53").Value = 162.57144
$ws.Range("N53").Value = -1544
$ws.Range("H62").Value = 6372.375
$ws.Range("J62").Value = 6500
$ws.Range("L62").Value = 6500
$ws.Range("N62").Value = -7748
$ws.Range("H65").Value = 6372.375
$ws.Range("J65").Value = 6500
$ws.Range("L65").Value = 32500
$ws.Range("N65").Value = -38740
$ws.Range("H98").Value = 53376
$ws.Range("I98").Value = 110778
$ws.Range("J98").Value = 12374.571
$ws.Range("K98").Value = 110778
$ws.Range("L98").Value = 12374.571
$ws.Range("M98").Value = -109280
$ws.Range("N98").Value = -15370.571
$ws.Range("H107").Value = 10992.6
$ws.Range("I107").Value = 11992.889
$ws.Range("K107").Value = 11992.889
$ws.Range("M107").Value = -10072.889
$ws.Range("H113").Value = 11665.637
$ws.Range("I113").Value = 12567.294
$ws.Range("K113").Value = 12567.294
$ws.Range("M113").Value = -9313.294
$ws.Range("H122").Value = 53376
$ws.Range("I122").Value = 110778
$ws.Range("J122").Value = 12374.571
$ws.Range("K122").Value = 332334
$ws.Range("L122").Value = 37123.713
$ws.Range("M122").Value = -329884
$ws.Range("N122").Value = -42023.713
$ws.Range("H123").Value = 113778.336
$ws.Range("J123").Value = 113778.336
$ws.Range("L123").Value = 113778.336
$ws.Range("N123").Value = -123578.336
$ws.Range("H125").Value = 3171
$ws.Range("I125").Value = 527.5
$ws.Range("J125").Value = 4933.3335
$ws.Range("K125").Value = 4747.5
$ws.Range("L125").Value = 44400.0015
$ws.Range("M125").Value = -2287.5
$ws.Range("N125").Value = -49320.0015
$ws.Range("H127").Value = 1186.6666
$ws.Range("I127").Value = 1186.6666
$ws.Range("K127").Value = 3559.9998
$ws.Range("M127").Value = 1400.0002
$ws.Range("H129").Value = 1837.4
$ws.Range("I129").Value = 1723.7273
$ws.Range("J129").Value = 2150
$ws.Range("K129").Value = 5171.1819
$ws.Range("L129").Value = 6450
$ws.Range("M129").Value = -171.1818999999996
$ws.Range("N129").Value = -16450

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1666.3334
$ws.Range("I4").Value = 1999
$ws.Range("J4").Value = 1500
$ws.Range("K4").Value = 1999
$ws.Range("L4").Value = 1500
$ws.Range("M4").Value = -1883
$ws.Range("N4").Value = -1732
$ws.Range("H37").Value = 30017
$ws.Range("H45").Value = 157088.64
$ws.Range("I45").Value = 217589.2
$ws.Range("K45").Value = 217589.2
$ws.Range("M45").Value = -217212.2
$ws.Range("H80").Value = 75823.25
$ws.Range("J80").Value = 75823.25
$ws.Range("L80").Value = 75823.25
$ws.Range("N80").Value = -77819.25
$ws.Range("H83").Value = 75823.25
$ws.Range("J83").Value = 75823.25
$ws.Range("L83").Value = 227469.75
$ws.Range("N83").Value = -237453.75
$ws.Range("H110").Value = 3100
$ws.Range("J110").Value = 3500
$ws.Range("L110").Value = 3500
$ws.Range("N110").Value = -7590
$ws.Range("H122").Value = 829088.1
$ws.Range("I122").Value = 4924.9
$ws.Range("K122").Value = 14774.7
$ws.Range("M122").Value = -12324.7
$ws.Range("H141").Value = 65000
$ws.Range("J141").Value = 65000
$ws.Range("L141").Value = 65000
$ws.Range("N141").Value = -75360

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H43").Value = 236167.9
$ws.Range("J43").Value = 236167.9
$ws.Range("L43").Value = 236167.9
$ws.Range("N43").Value = -236529.9
$ws.Range("H86").Value = 4999.923
$ws.Range("J86").Value = 2459.182
$ws.Range("L86").Value = 2459.182
$ws.Range("N86").Value = -4705.182
$ws.Range("H89").Value = 4999.923
$ws.Range("J89").Value = 2459.182
$ws.Range("L89").Value = 12295.91
$ws.Range("N89").Value = -23527.91
$ws.Range("H134").Value = 5007.778
$ws.Range("I134").Value = 5045.1465
$ws.Range("K134").Value = 15135.4395
$ws.Range("M134").Value = -12600.4395

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 9092001
$ws.Range("I7").Value = 1440.25
$ws.Range("J7").Value = 33333496
$ws.Range("K7").Value = 1440.25
$ws.Range("L7").Value = 33333496
$ws.Range("M7").Value = -1327.25
$ws.Range("N7").Value = -33333722
$ws.Range("H16").Value = 2631.7
$ws.Range("I16").Value = 3371.3333
$ws.Range("J16").Value = 1522.25
$ws.Range("K16").Value = 3371.3333
$ws.Range("L16").Value = 1522.25
$ws.Range("M16").Value = -3084.3333
$ws.Range("N16").Value = -2096.25
$ws.Range("H113").Value = 2631.7
$ws.Range("I113").Value = 3371.3333
$ws.Range("J113").Value = 1522.25
$ws.Range("K113").Value = 3371.3333
$ws.Range("L113").Value = 1522.25
$ws.Range("M113").Value = -1201.3333
$ws.Range("N113").Value = -5862.25
$ws.Range("H119").Value = 49999
$ws.Range("J119").Value = 49997
$ws.Range("L119").Value = 49997
$ws.Range("N119").Value = -59673
$ws.Range("H120").Value = 51597.715
$ws.Range("I120").Value = 50296
$ws.Range("J120").Value = 53333.332
$ws.Range("K120").Value = 50296
$ws.Range("L120").Value = 53333.332
$ws.Range("M120").Value = -46667
$ws.Range("N120").Value = -60591.332
$ws.Range("H132").Value = 27169.521
$ws.Range("I132").Value = 10894.25
$ws.Range("K132").Value = 32682.75
$ws.Range("M132").Value = -30152.75

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 299
$ws.Range("I32").Value = 299
$ws.Range("K32").Value = 897
$ws.Range("M32").Value = -614
$ws.Range("H68").Value = 5670.512
$ws.Range("J68").Value = 6581.9062
$ws.Range("L68").Value = 19745.7186
$ws.Range("N68").Value = -21367.7186
$ws.Range("H71").Value = 5670.512
$ws.Range("J71").Value = 6581.9062
$ws.Range("L71").Value = 59237.1558
$ws.Range("N71").Value = -67349.15580000001
$ws.Range("H107").Value = 2519.524
$ws.Range("I107").Value = 716.2857
$ws.Range("J107").Value = 3421.1428
$ws.Range("K107").Value = 2148.8571
$ws.Range("L107").Value = 10263.4284
$ws.Range("M107").Value = -228.8571000000002
$ws.Range("N107").Value = -14103.4284

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 52499.5
$ws.Range("I97").Value = 68333
$ws.Range("J97").Value = 4999
$ws.Range("K97").Value = 68333
$ws.Range("L97").Value = 4999
$ws.Range("M97").Value = -67837
$ws.Range("N97").Value = -5991
$ws.Range("H113").Value = 2449.25
$ws.Range("J113").Value = 2449.25
$ws.Range("L113").Value = 2449.25
$ws.Range("N113").Value = -6789.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1005.5
$ws.Range("I16").Value = 1017.2759
$ws.Range("J16").Value = 937.2
$ws.Range("K16").Value = 1017.2759
$ws.Range("L16").Value = 937.2
$ws.Range("M16").Value = -847.2759
$ws.Range("N16").Value = -1277.2
$ws.Range("H46").Value = 2460.1724
$ws.Range("I46").Value = 659.0833
$ws.Range("J46").Value = 3731.5293
$ws.Range("K46").Value = 659.0833
$ws.Range("L46").Value = 3731.5293
$ws.Range("M46").Value = -471.0833
$ws.Range("N46").Value = -4107.5293
$ws.Range("H61").Value = 3190.25
$ws.Range("I61").Value = 3421.2
$ws.Range("J61").Value = 3025.2856
$ws.Range("K61").Value = 3421.2
$ws.Range("L61").Value = 3025.2856
$ws.Range("M61").Value = -3219.2
$ws.Range("N61").Value = -3429.2856
$ws.Range("H93").Value = 6647.8887
$ws.Range("I93").Value = 6647.8887
$ws.Range("K93").Value = 6647.8887
$ws.Range("M93").Value = -5399.8887
$ws.Range("H113").Value = 3190.25
$ws.Range("I113").Value = 3421.2
$ws.Range("J113").Value = 3025.2856
$ws.Range("K113").Value = 3421.2
$ws.Range("L113").Value = 3025.2856
$ws.Range("M113").Value = -1251.2
$ws.Range("N113").Value = -7365.2856
$ws.Range("H122").Value = 7366.6665
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 31095.715
$ws.Range("I81").Value = 41349
$ws.Range("J81").Value = 5462.5
$ws.Range("K81").Value = 82698
$ws.Range("L81").Value = 10925
$ws.Range("M81").Value = -81637
$ws.Range("N81").Value = -13047
$ws.Range("H84").Value = 31095.715
$ws.Range("I84").Value = 41349
$ws.Range("J84").Value = 5462.5
$ws.Range("K84").Value = 413490
$ws.Range("L84").Value = 54625
$ws.Range("M84").Value = -408186
$ws.Range("N84").Value = -65233
$ws.Range("H107").Value = 3757.2222
$ws.Range("I107").Value = 4817.5
$ws.Range("K107").Value = 14452.5
$ws.Range("M107").Value = -12532.5
$ws.Range("H122").Value = 10997.5
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 10997.5
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 32992.5
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -37892.5

Write-Host "All edits applied"